$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the import line.
# ---------------------------------------------------------------------------
$d.Paragraphs(1).Range.Text = 'import re'

# ---------------------------------------------------------------------------
# 2) Delete the old unittest-based test class block (paragraphs 2-30, which
#    sit between "import unittest"/"import re" and the blank line that
#    precedes "def validate_password(password):").
# ---------------------------------------------------------------------------
$startRange = $d.Paragraphs(2).Range.Start
$endRange = $d.Paragraphs(30).Range.End
$d.Range($startRange, $endRange).Delete()

# After the delete, the paragraph layout is:
#   1  import re
#   2  (blank)
#   3  def validate_password(password):
#   4      if len(password) < 8:
#   5          return 'Rejected - Must be at least 8 characters'
#   6      (4 spaces)
#   7      if not any(char.isdigit() for char in password):
#   8          return 'Rejected - Must include at least one number'
#   9  (blank)
#  10      if not any(char in '!@#$%^&*' for char in password):
#  11          return 'Rejected - Must include at least one special character'
#  12  (blank)
#  13      return 'Accepted'
#  14  (blank)
#  15  if __name__ == '__main__':

# ---------------------------------------------------------------------------
# 3) Rewrite the body of validate_password().
# ---------------------------------------------------------------------------
$d.Paragraphs(4).Range.Text = '  if len(password) < 8:'
$d.Paragraphs(5).Range.Text = '    return False'
$d.Paragraphs(6).Range.Text = '  '
$d.Paragraphs(7).Range.Text = '  if not re.search(r''\d'', password):'
$d.Paragraphs(8).Range.Text = '    return False'
$d.Paragraphs(10).Range.Text = '  if not re.search(r''[!@#$%^&*(),.?":{}|<>]'', password):'
$d.Paragraphs(11).Range.Text = '    return False'
$d.Paragraphs(13).Range.Text = '  return True'

# ---------------------------------------------------------------------------
# 4) Append the new test_passwords() block after "return True" (paragraph 13),
#    replacing the old trailing blank line + "if __name__" footer, which get
#    re-appended at the very end of the inserted block.
# ---------------------------------------------------------------------------
$newLines = @(
    '',
    '',
    'def test_passwords():',
    '  ',
    '  test1 = validate_password(''abc1$'')',
    '  assert test1 == False',
    '',
    '  test2 = validate_password(''abcd@xyz'') ',
    '  assert test2 == False',
    '',
    '  test3 = validate_password(''abcd1234'')',
    '  assert test3 == False',
    '',
    '  test4 = validate_password(''abcd@xyz'')',
    '  assert test4 == False',
    '',
    '  test5 = validate_password(''abc1@def'')',
    '  assert test5 == True',
    '',
    '  test6 = validate_password(''MyPass123!'')',
    '  assert test6 == True',
    '',
    '  test7 = validate_password(''1234@5678'')',
    '  assert test7 == True',
    '',
    '  test8 = validate_password(''abcdefgh'')',
    '  assert test8 == False',
    '',
    '  test9 = validate_password(''abcd1234'')',
    '  assert test9 == False',
    '',
    '  test10 = validate_password(''Ab1$xyz9'')',
    '  assert test10 == True',
    ''
)

$p = $d.Paragraphs(13)
foreach ($line in $newLines) {
    $p.Range.InsertParagraphAfter()
    $p = $p.Next()
    $p.Range.Text = $line
}

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
